$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# remain text (matching the source data, which stores all values as
# inline/shared strings, not numbers).
$forceTextCells = @("D5,D6,D10,D11,D12,D15,D17,D19,D23,D25,D27,D28,D29,D30,D31,D32,D33,D34,D35,D38,D39,D40,D41,D42,D45,D46,D48,D49,D50,D51".Split(","))
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "44.150.11"

# Row 3
$ws.Range("D3").Value = "2.257.04"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "308.47"
$ws.Range("E5").Value = "  -4.27%  "

# Row 6
$ws.Range("D6").Value = "98.98"
$ws.Range("E6").Value = "  -2.50%  "

# Row 7
$ws.Range("E7").Value = "  -0.67%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("E9").Value = "  -3.29%  "

# Row 10
$ws.Range("D10").Value = "35.68"
$ws.Range("E10").Value = "  -4.43%  "

# Row 11
$ws.Range("D11").Value = "0.0823"
$ws.Range("E11").Value = "  -0.82%  "

# Row 12
$ws.Range("D12").Value = "7.35"
$ws.Range("E12").Value = "  -4.34%  "

# Row 13
$ws.Range("E13").Value = "  -1.88%  "

# Row 14
$ws.Range("D14").Value = "2.601.02"
$ws.Range("E14").Value = "  +0.56%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.842"
$ws.Range("E15").Value = "  -1.58%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.252.18"
$ws.Range("E16").Value = "  +0.35%  "

# Row 17
$ws.Range("D17").Value = "13.91"
$ws.Range("E17").Value = "  -1.85%  "

# Row 18
$ws.Range("D18").Value = "44.059.00"
$ws.Range("E18").Value = "  +0.99%  "

# Row 19
$ws.Range("D19").Value = "12.88"
$ws.Range("E19").Value = "  -6.35%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0976"
$ws.Range("E20").Value = "  -0.84%  "

# Row 21
$ws.Range("E21").Value = "  -1.19%  "

# Row 22
$ws.Range("E22").Value = "  +0.39%  "

# Row 23
$ws.Range("D23").Value = "241.28"
$ws.Range("E23").Value = "  +2.02%  "

# Row 24
$ws.Range("E24").Value = "  -6.46%  "

# Row 25
$ws.Range("D25").Value = "1.98"
$ws.Range("E25").Value = "  -8.01%  "

# Row 26
$ws.Range("E26").Value = "  +0.36%  "

# Row 27
$ws.Range("D27").Value = "10.18"
$ws.Range("E27").Value = "  +1.05%  "

# Row 28
$ws.Range("D28").Value = "37.68"
$ws.Range("E28").Value = "  +1.65%  "

# Row 29
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  -0.39%  "

# Row 30
$ws.Range("D30").Value = "6.19"
$ws.Range("E30").Value = "  -1.48%  "

# Row 31
$ws.Range("D31").Value = "20.15"
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("D32").Value = "157.21"
$ws.Range("E32").Value = "  -1.72%  "

# Row 33
$ws.Range("D33").Value = "3.57"
$ws.Range("E33").Value = "  +12.10%  "

# Row 34
$ws.Range("D34").Value = "0.0824"
$ws.Range("E34").Value = "  -3.17%  "

# Row 35
$ws.Range("D35").Value = "2.66"
$ws.Range("E35").Value = "  -1.10%  "

# Row 36
$ws.Range("E36").Value = "  -0.30%  "

# Row 37
$ws.Range("E37").Value = "  -5.55%  "

# Row 38
$ws.Range("D38").Value = "1.87"
$ws.Range("E38").Value = "  -3.04%  "

# Row 39
$ws.Range("D39").Value = "15.80"
$ws.Range("E39").Value = "  +2.85%  "

# Row 40
$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  -8.86%  "

# Row 41
$ws.Range("D41").Value = "3.40"
$ws.Range("E41").Value = "  -10.37%  "

# Row 42
$ws.Range("D42").Value = "0.0306"
$ws.Range("E42").Value = "  -3.61%  "

# Row 43
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("D44").Value = "1.776.18"
$ws.Range("E44").Value = "  -1.93%  "

# Row 45
$ws.Range("D45").Value = "88.15"
$ws.Range("E45").Value = "  +6.53%  "

# Row 46
$ws.Range("D46").Value = "0.193"
$ws.Range("E46").Value = "  -3.63%  "

# Row 47
$ws.Range("E47").Value = "  -1.54%  "

# Row 48
$ws.Range("D48").Value = "101.83"
$ws.Range("E48").Value = "  -1.91%  "

# Row 49
$ws.Range("D49").Value = "8.29"
$ws.Range("E49").Value = "  -2.03%  "

# Row 50
$ws.Range("D50").Value = "70.39"
$ws.Range("E50").Value = "  -5.53%  "

# Row 51
$ws.Range("D51").Value = "55.62"
$ws.Range("E51").Value = "  -5.50%  "

# Restore default style on forced-text cells (NumberFormat="@" assigns a
# new style index; reset back to Normal so no visible style changes linger).
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
